$wb = $excel.ActiveWorkbook

# This script applies a refreshed market-price/profit data snapshot
# (currentAveragePrice*, Leve*Price*, LeveProfit* columns) to several
# worksheets, mirroring a scheduled data-update run.

$ws = $wb.Worksheets.Item("ALC")
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
# Row 68
$ws.Range("H68").Value = 69999
$ws.Range("J68").Value = 69999
$ws.Range("L68").Value = 69999
$ws.Range("N68").Value = -71497
# Row 71
$ws.Range("H71").Value = 69999
$ws.Range("J71").Value = 69999
$ws.Range("L71").Value = 209997
$ws.Range("N71").Value = -217485
# Row 75
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872
# Row 78
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360
# Row 80
$ws.Range("H80").Value = 1708.1538
$ws.Range("I80").Value = 2663.8572
$ws.Range("J80").Value = 593.1667
$ws.Range("K80").Value = 7991.571599999999
$ws.Range("L80").Value = 1779.5001
$ws.Range("M80").Value = -6993.571599999999
$ws.Range("N80").Value = -3775.5001
# Row 83
$ws.Range("H83").Value = 1708.1538
$ws.Range("I83").Value = 2663.8572
$ws.Range("J83").Value = 593.1667
$ws.Range("K83").Value = 23974.7148
$ws.Range("L83").Value = 5338.5003
$ws.Range("M83").Value = -18982.7148
$ws.Range("N83").Value = -15322.5003
# Row 86
$ws.Range("H86").Value = 4249.1333
$ws.Range("I86").Value = 6514.625
$ws.Range("J86").Value = 1660
$ws.Range("K86").Value = 6514.625
$ws.Range("L86").Value = 1660
$ws.Range("M86").Value = -5391.625
$ws.Range("N86").Value = -3906
# Row 88
$ws.Range("H88").Value = 23079782
$ws.Range("J88").Value = 2991597.8
$ws.Range("L88").Value = 2991597.8
$ws.Range("N88").Value = -2992409.8
# Row 89
$ws.Range("H89").Value = 4249.1333
$ws.Range("I89").Value = 6514.625
$ws.Range("J89").Value = 1660
$ws.Range("K89").Value = 32573.125
$ws.Range("L89").Value = 8300
$ws.Range("M89").Value = -26957.125
$ws.Range("N89").Value = -19532
# Row 91
$ws.Range("H91").Value = 23079782
$ws.Range("J91").Value = 2991597.8
$ws.Range("L91").Value = 2991597.8
$ws.Range("N91").Value = -2994405.8
# Row 92
$ws.Range("H92").Value = 847.1
$ws.Range("I92").Value = 821.94116
$ws.Range("K92").Value = 821.94116
$ws.Range("M92").Value = 426.05884
# Row 98
$ws.Range("H98").Value = 1205.258
$ws.Range("I98").Value = 912.1
$ws.Range("K98").Value = 912.1
$ws.Range("M98").Value = 585.9
# Row 99
$ws.Range("H99").Value = 2678.2727
$ws.Range("I99").Value = 245.375
$ws.Range("K99").Value = 736.125
$ws.Range("M99").Value = 761.875
# Row 122
$ws.Range("H122").Value = 1205.258
$ws.Range("I122").Value = 912.1
$ws.Range("K122").Value = 2736.3
$ws.Range("M122").Value = -286.3000000000002

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 2132.1365
$ws.Range("J88").Value = 2536.6365
$ws.Range("L88").Value = 2536.6365
$ws.Range("N88").Value = -3348.6365
# Row 91
$ws.Range("H91").Value = 2132.1365
$ws.Range("J91").Value = 2536.6365
$ws.Range("L91").Value = 2536.6365
$ws.Range("N91").Value = -5344.636500000001

$ws = $wb.Worksheets.Item("CRP")
# Row 69
$ws.Range("H69").Value = 16500
$ws.Range("I69").Value = 16500
$ws.Range("K69").Value = 16500
$ws.Range("M69").Value = -15751
# Row 72
$ws.Range("H72").Value = 16500
$ws.Range("I72").Value = 16500
$ws.Range("K72").Value = 49500
$ws.Range("M72").Value = -45756
# Row 94
$ws.Range("H94").Value = 1689.1305
$ws.Range("I94").Value = 1748.5385
$ws.Range("J94").Value = 1611.9
$ws.Range("K94").Value = 1748.5385
$ws.Range("L94").Value = 1611.9
$ws.Range("M94").Value = -1297.5385
$ws.Range("N94").Value = -2513.9

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1839.2727
$ws.Range("I131").Value = 1309.75
$ws.Range("K131").Value = 3929.25
$ws.Range("M131").Value = 1110.75
# Row 132
$ws.Range("H132").Value = 1045.8334
$ws.Range("I132").Value = 785
$ws.Range("J132").Value = 1176.25
$ws.Range("K132").Value = 7065
$ws.Range("L132").Value = 10586.25
$ws.Range("M132").Value = -4535
$ws.Range("N132").Value = -15646.25

$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 1999.5
$ws.Range("I14").Value = 1999.5
$ws.Range("K14").Value = 1999.5
$ws.Range("M14").Value = -1827.5
# Row 88
$ws.Range("H88").Value = 21333.334
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 21333.334
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 21333.334
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -22189.334
# Row 91
$ws.Range("H91").Value = 21333.334
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 21333.334
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 21333.334
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -24297.334

$ws = $wb.Worksheets.Item("WVR")
# Row 68
$ws.Range("H68").Value = 271
$ws.Range("J68").Value = 271
$ws.Range("L68").Value = 271
$ws.Range("N68").Value = -1893
# Row 71
$ws.Range("H71").Value = 271
$ws.Range("J71").Value = 271
$ws.Range("L71").Value = 813
$ws.Range("N71").Value = -8925
# Row 75
$ws.Range("H75").Value = 109449.75
$ws.Range("I75").Value = 37800
$ws.Range("K75").Value = 37800
$ws.Range("M75").Value = -36864
# Row 78
$ws.Range("H78").Value = 109449.75
$ws.Range("I78").Value = 37800
$ws.Range("K78").Value = 113400
$ws.Range("M78").Value = -108720
# Row 80
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
# Row 81
$ws.Range("H81").Value = 5199.4
$ws.Range("I81").Value = 3499.25
$ws.Range("J81").Value = 12000
$ws.Range("K81").Value = 6998.5
$ws.Range("L81").Value = 24000
$ws.Range("M81").Value = -5937.5
$ws.Range("N81").Value = -26122
# Row 83
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
# Row 84
$ws.Range("H84").Value = 5199.4
$ws.Range("I84").Value = 3499.25
$ws.Range("J84").Value = 12000
$ws.Range("K84").Value = 34992.5
$ws.Range("L84").Value = 120000
$ws.Range("M84").Value = -29688.5
$ws.Range("N84").Value = -130608
# Row 126
$ws.Range("H126").Value = 2565.1025
$ws.Range("I126").Value = 2644.7666
$ws.Range("K126").Value = 7934.2998
$ws.Range("M126").Value = -5464.2998
# Row 132
$ws.Range("H132").Value = 9809148
$ws.Range("I132").Value = 11365865
$ws.Range("J132").Value = 24071.143
$ws.Range("K132").Value = 34097595
$ws.Range("L132").Value = 72213.429
$ws.Range("M132").Value = -34095065
$ws.Range("N132").Value = -77273.429
